# Update the GSC "Breadcrumbs" export: roll the rolling date window
# forward by one day (drop 2025-11-18, add 2026-02-13 / 2026-02-14) and
# refresh the per-day "Valid" counters; "Invalid" stays 0 throughout.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column A holds plain-text yyyy-MM-dd labels (not real dates). Force
# text formatting first so Excel does not auto-convert the strings we
# assign below into date serials.
$ws.Range("A2:A89").NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = "2025-11-19"
$ws.Cells.Item(2, 3).Value = 26.0
$ws.Cells.Item(3, 1).Value = "2025-11-20"
$ws.Cells.Item(3, 3).Value = 25.0
$ws.Cells.Item(4, 1).Value = "2025-11-21"
$ws.Cells.Item(4, 3).Value = 26.0
$ws.Cells.Item(5, 1).Value = "2025-11-22"
$ws.Cells.Item(5, 3).Value = 26.0
$ws.Cells.Item(6, 1).Value = "2025-11-23"
$ws.Cells.Item(6, 3).Value = 25.0
$ws.Cells.Item(7, 1).Value = "2025-11-24"
$ws.Cells.Item(7, 3).Value = 25.0
$ws.Cells.Item(8, 1).Value = "2025-11-25"
$ws.Cells.Item(8, 3).Value = 27.0
$ws.Cells.Item(9, 1).Value = "2025-11-26"
$ws.Cells.Item(9, 3).Value = 27.0
$ws.Cells.Item(10, 1).Value = "2025-11-27"
$ws.Cells.Item(10, 3).Value = 27.0
$ws.Cells.Item(11, 1).Value = "2025-11-28"
$ws.Cells.Item(11, 3).Value = 27.0
$ws.Cells.Item(12, 1).Value = "2025-11-29"
$ws.Cells.Item(12, 3).Value = 27.0
$ws.Cells.Item(13, 1).Value = "2025-11-30"
$ws.Cells.Item(13, 3).Value = 27.0
$ws.Cells.Item(14, 1).Value = "2025-12-01"
$ws.Cells.Item(14, 3).Value = 27.0
$ws.Cells.Item(15, 1).Value = "2025-12-02"
$ws.Cells.Item(15, 3).Value = 27.0
$ws.Cells.Item(16, 1).Value = "2025-12-03"
$ws.Cells.Item(16, 3).Value = 26.0
$ws.Cells.Item(17, 1).Value = "2025-12-04"
$ws.Cells.Item(17, 3).Value = 25.0
$ws.Cells.Item(18, 1).Value = "2025-12-05"
$ws.Cells.Item(18, 3).Value = 25.0
$ws.Cells.Item(19, 1).Value = "2025-12-06"
$ws.Cells.Item(19, 3).Value = 25.0
$ws.Cells.Item(20, 1).Value = "2025-12-07"
$ws.Cells.Item(20, 3).Value = 26.0
$ws.Cells.Item(21, 1).Value = "2025-12-08"
$ws.Cells.Item(21, 3).Value = 26.0
$ws.Cells.Item(22, 1).Value = "2025-12-09"
$ws.Cells.Item(22, 3).Value = 27.0
$ws.Cells.Item(23, 1).Value = "2025-12-10"
$ws.Cells.Item(23, 3).Value = 29.0
$ws.Cells.Item(24, 1).Value = "2025-12-11"
$ws.Cells.Item(24, 3).Value = 29.0
$ws.Cells.Item(25, 1).Value = "2025-12-12"
$ws.Cells.Item(25, 3).Value = 30.0
$ws.Cells.Item(26, 1).Value = "2025-12-13"
$ws.Cells.Item(26, 3).Value = 30.0
$ws.Cells.Item(27, 1).Value = "2025-12-14"
$ws.Cells.Item(27, 3).Value = 31.0
$ws.Cells.Item(28, 1).Value = "2025-12-15"
$ws.Cells.Item(28, 3).Value = 31.0
$ws.Cells.Item(29, 1).Value = "2025-12-16"
$ws.Cells.Item(29, 3).Value = 32.0
$ws.Cells.Item(30, 1).Value = "2025-12-17"
$ws.Cells.Item(30, 3).Value = 31.0
$ws.Cells.Item(31, 1).Value = "2025-12-18"
$ws.Cells.Item(31, 3).Value = 31.0
$ws.Cells.Item(32, 1).Value = "2025-12-19"
$ws.Cells.Item(32, 3).Value = 32.0
$ws.Cells.Item(33, 1).Value = "2025-12-20"
$ws.Cells.Item(33, 3).Value = 32.0
$ws.Cells.Item(34, 1).Value = "2025-12-21"
$ws.Cells.Item(34, 3).Value = 32.0
$ws.Cells.Item(35, 1).Value = "2025-12-22"
$ws.Cells.Item(35, 3).Value = 32.0
$ws.Cells.Item(36, 1).Value = "2025-12-23"
$ws.Cells.Item(36, 3).Value = 30.0
$ws.Cells.Item(37, 1).Value = "2025-12-24"
$ws.Cells.Item(37, 3).Value = 31.0
$ws.Cells.Item(38, 1).Value = "2025-12-25"
$ws.Cells.Item(38, 3).Value = 32.0
$ws.Cells.Item(39, 1).Value = "2025-12-26"
$ws.Cells.Item(39, 3).Value = 32.0
$ws.Cells.Item(40, 1).Value = "2025-12-27"
$ws.Cells.Item(40, 3).Value = 28.0
$ws.Cells.Item(41, 1).Value = "2025-12-28"
$ws.Cells.Item(41, 3).Value = 28.0
$ws.Cells.Item(42, 1).Value = "2025-12-29"
$ws.Cells.Item(42, 3).Value = 28.0
$ws.Cells.Item(43, 1).Value = "2025-12-30"
$ws.Cells.Item(43, 3).Value = 28.0
$ws.Cells.Item(44, 1).Value = "2025-12-31"
$ws.Cells.Item(44, 3).Value = 30.0
$ws.Cells.Item(45, 1).Value = "2026-01-01"
$ws.Cells.Item(45, 3).Value = 29.0
$ws.Cells.Item(46, 1).Value = "2026-01-02"
$ws.Cells.Item(46, 3).Value = 28.0
$ws.Cells.Item(47, 1).Value = "2026-01-03"
$ws.Cells.Item(47, 3).Value = 28.0
$ws.Cells.Item(48, 1).Value = "2026-01-04"
$ws.Cells.Item(48, 3).Value = 27.0
$ws.Cells.Item(49, 1).Value = "2026-01-05"
$ws.Cells.Item(49, 3).Value = 27.0
$ws.Cells.Item(50, 1).Value = "2026-01-06"
$ws.Cells.Item(50, 3).Value = 27.0
$ws.Cells.Item(51, 1).Value = "2026-01-07"
$ws.Cells.Item(51, 3).Value = 27.0
$ws.Cells.Item(52, 1).Value = "2026-01-08"
$ws.Cells.Item(52, 3).Value = 27.0
$ws.Cells.Item(53, 1).Value = "2026-01-09"
$ws.Cells.Item(53, 3).Value = 27.0
$ws.Cells.Item(54, 1).Value = "2026-01-10"
$ws.Cells.Item(54, 3).Value = 26.0
$ws.Cells.Item(55, 1).Value = "2026-01-11"
$ws.Cells.Item(55, 3).Value = 26.0
$ws.Cells.Item(56, 1).Value = "2026-01-12"
$ws.Cells.Item(56, 3).Value = 26.0
$ws.Cells.Item(57, 1).Value = "2026-01-13"
$ws.Cells.Item(57, 3).Value = 26.0
$ws.Cells.Item(58, 1).Value = "2026-01-14"
$ws.Cells.Item(58, 3).Value = 26.0
$ws.Cells.Item(59, 1).Value = "2026-01-15"
$ws.Cells.Item(59, 3).Value = 26.0
$ws.Cells.Item(60, 1).Value = "2026-01-16"
$ws.Cells.Item(60, 3).Value = 26.0
$ws.Cells.Item(61, 1).Value = "2026-01-17"
$ws.Cells.Item(61, 3).Value = 25.0
$ws.Cells.Item(62, 1).Value = "2026-01-18"
$ws.Cells.Item(62, 3).Value = 25.0
$ws.Cells.Item(63, 1).Value = "2026-01-19"
$ws.Cells.Item(63, 3).Value = 26.0
$ws.Cells.Item(64, 1).Value = "2026-01-20"
$ws.Cells.Item(64, 3).Value = 25.0
$ws.Cells.Item(65, 1).Value = "2026-01-21"
$ws.Cells.Item(65, 3).Value = 24.0
$ws.Cells.Item(66, 1).Value = "2026-01-22"
$ws.Cells.Item(66, 3).Value = 23.0
$ws.Cells.Item(67, 1).Value = "2026-01-23"
$ws.Cells.Item(67, 3).Value = 24.0
$ws.Cells.Item(68, 1).Value = "2026-01-24"
$ws.Cells.Item(68, 3).Value = 24.0
$ws.Cells.Item(69, 1).Value = "2026-01-25"
$ws.Cells.Item(69, 3).Value = 24.0
$ws.Cells.Item(70, 1).Value = "2026-01-26"
$ws.Cells.Item(70, 3).Value = 25.0
$ws.Cells.Item(71, 1).Value = "2026-01-27"
$ws.Cells.Item(71, 3).Value = 26.0
$ws.Cells.Item(72, 1).Value = "2026-01-28"
$ws.Cells.Item(72, 3).Value = 26.0
$ws.Cells.Item(73, 1).Value = "2026-01-29"
$ws.Cells.Item(73, 3).Value = 28.0
$ws.Cells.Item(74, 1).Value = "2026-01-30"
$ws.Cells.Item(74, 3).Value = 28.0
$ws.Cells.Item(75, 1).Value = "2026-01-31"
$ws.Cells.Item(75, 3).Value = 28.0
$ws.Cells.Item(76, 1).Value = "2026-02-01"
$ws.Cells.Item(76, 3).Value = 28.0
$ws.Cells.Item(77, 1).Value = "2026-02-02"
$ws.Cells.Item(77, 3).Value = 28.0
$ws.Cells.Item(78, 1).Value = "2026-02-03"
$ws.Cells.Item(78, 3).Value = 28.0
$ws.Cells.Item(79, 1).Value = "2026-02-04"
$ws.Cells.Item(79, 3).Value = 28.0
$ws.Cells.Item(80, 1).Value = "2026-02-05"
$ws.Cells.Item(80, 3).Value = 28.0
$ws.Cells.Item(81, 1).Value = "2026-02-06"
$ws.Cells.Item(81, 3).Value = 28.0
$ws.Cells.Item(82, 1).Value = "2026-02-07"
$ws.Cells.Item(82, 3).Value = 27.0
$ws.Cells.Item(83, 1).Value = "2026-02-08"
$ws.Cells.Item(83, 3).Value = 28.0
$ws.Cells.Item(84, 1).Value = "2026-02-09"
$ws.Cells.Item(84, 3).Value = 28.0
$ws.Cells.Item(85, 1).Value = "2026-02-10"
$ws.Cells.Item(85, 3).Value = 29.0
$ws.Cells.Item(86, 1).Value = "2026-02-11"
$ws.Cells.Item(86, 3).Value = 30.0
$ws.Cells.Item(87, 1).Value = "2026-02-12"
$ws.Cells.Item(87, 3).Value = 30.0
$ws.Cells.Item(88, 1).Value = "2026-02-13"
$ws.Cells.Item(88, 3).Value = 30.0
$ws.Cells.Item(89, 1).Value = "2026-02-14"
$ws.Cells.Item(89, 3).Value = 31.0

# New row 89 (2026-02-14) also needs its "Invalid" counter (always 0).
$ws.Cells.Item(89, 2).Value = 0
